$wb = $excel.ActiveWorkbook

$wsWorkers = $wb.Worksheets.Item("Workers")
$wsBios = $wb.Worksheets.Item("Bios")
$wsSkills = $wb.Worksheets.Item("Skills")
$wsContracts = $wb.Worksheets.Item("Contracts")
$wsNotes = $wb.Worksheets.Item("Notes")

# --- Workers: new worker "Mr Bubba" (UID 2964), row 5 ---
$wsWorkers.Cells.Item(5, 1).Value = 2964
$wsWorkers.Cells.Item(5, 2).Value = 0
$wsWorkers.Cells.Item(5, 3).Value = 0
$wsWorkers.Cells.Item(5, 4).Value = -1
$wsWorkers.Cells.Item(5, 5).Value = "Mr Bubba"
$wsWorkers.Cells.Item(5, 6).Value = "Mr"
$wsWorkers.Cells.Item(5, 7).Value = 1
$wsWorkers.Cells.Item(5, 8).Value = 1
$wsWorkers.Cells.Item(5, 9).Value = 1
$wsWorkers.Cells.Item(5, 10).Value = 2
$wsWorkers.Cells.Item(5, 11).Value = 0
$wsWorkers.Cells.Item(5, 12).Value = 32116
$wsWorkers.Cells.Item(5, 12).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWorkers.Cells.Item(5, 13).Value = 39448
$wsWorkers.Cells.Item(5, 13).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWorkers.Cells.Item(5, 14).NumberFormat = "@"
$wsWorkers.Cells.Item(5, 14).Value = "1666-01-01"
$wsWorkers.Cells.Item(5, 14).ClearFormats()
$wsWorkers.Cells.Item(5, 15).Value = 1
$wsWorkers.Cells.Item(5, 16).Value = 25
$wsWorkers.Cells.Item(5, 17).Value = 262
$wsWorkers.Cells.Item(5, 18).Value = 242
$wsWorkers.Cells.Item(5, 19).Value = 296
$wsWorkers.Cells.Item(5, 20).Value = "mrbubba.jpg"
$wsWorkers.Cells.Item(5, 21).Value = 1
$wsWorkers.Cells.Item(5, 22).Value = 9
$wsWorkers.Cells.Item(5, 23).Value = 1
$wsWorkers.Cells.Item(5, 24).Value = 0
$wsWorkers.Cells.Item(5, 25).Value = 0
$wsWorkers.Cells.Item(5, 26).Value = 0
$wsWorkers.Cells.Item(5, 27).Value = 0
$wsWorkers.Cells.Item(5, 28).Value = 0
$wsWorkers.Cells.Item(5, 29).Value = 6
$wsWorkers.Cells.Item(5, 30).Value = 0
$wsWorkers.Cells.Item(5, 31).Value = 0
$wsWorkers.Cells.Item(5, 32).Value = 0
$wsWorkers.Cells.Item(5, 33).Value = -1
$wsWorkers.Cells.Item(5, 34).Value = -1
$wsWorkers.Cells.Item(5, 35).Value = -1
$wsWorkers.Cells.Item(5, 36).Value = -1
$wsWorkers.Cells.Item(5, 37).Value = -1
$wsWorkers.Cells.Item(5, 38).Value = -1
$wsWorkers.Cells.Item(5, 39).Value = -1
$wsWorkers.Cells.Item(5, 40).Value = -1
$wsWorkers.Cells.Item(5, 41).Value = 4
$wsWorkers.Cells.Item(5, 42).Value = 4
$wsWorkers.Cells.Item(5, 43).Value = 4
$wsWorkers.Cells.Item(5, 44).Value = 4
$wsWorkers.Cells.Item(5, 45).Value = 4
$wsWorkers.Cells.Item(5, 46).Value = 4
$wsWorkers.Cells.Item(5, 47).Value = 4
$wsWorkers.Cells.Item(5, 48).Value = 1
$wsWorkers.Cells.Item(5, 49).Value = 0
$wsWorkers.Cells.Item(5, 50).Value = -1
$wsWorkers.Cells.Item(5, 51).Value = 0
$wsWorkers.Cells.Item(5, 52).Value = 0
$wsWorkers.Cells.Item(5, 53).Value = 0
$wsWorkers.Cells.Item(5, 54).Value = 0
$wsWorkers.Cells.Item(5, 55).Value = 0
$wsWorkers.Cells.Item(5, 56).Value = 0
$wsWorkers.Cells.Item(5, 57).Value = 0
$wsWorkers.Cells.Item(5, 58).Value = 0
$wsWorkers.Cells.Item(5, 59).Value = 0
$wsWorkers.Cells.Item(5, 60).Value = 0
$wsWorkers.Cells.Item(5, 61).Value = 0
$wsWorkers.Cells.Item(5, 62).Value = 0
$wsWorkers.Cells.Item(5, 63).Value = -1
$wsWorkers.Cells.Item(5, 64).Value = "Bunny Brawler"
$wsWorkers.Cells.Item(5, 65).Value = 1
$wsWorkers.Cells.Item(5, 66).Value = "Bunny Villain"
$wsWorkers.Cells.Item(5, 67).Value = 1
$wsWorkers.Cells.Item(5, 68).Value = 0

# --- Bios: bio text for Mr Bubba, row 5 ---
$wsBios.Cells.Item(5, 1).Value = 2964
$wsBios.Cells.Item(5, 2).Value = "Introducing the enigmatic and mysterious professional wrestler known as Mr. Bubba. Hailing from parts unknown, Mr. Bubba is a male competitor who has taken the wrestling world by storm with his unique and captivating persona.
Sporting a bunny mask that covers his face, Mr. Bubba brings a sense of intrigue and unpredictability to the ring. His wrestling style is unlike anything fans have seen before, as he seamlessly blends elements of interpretive dance with traditional wrestling techniques, creating a one-of-a-kind performance that leaves audiences in awe.
Despite his unconventional approach, Mr. Bubba has quickly earned a reputation as a formidable competitor in the squared circle. His agility and athleticism are unmatched, allowing him to execute breathtaking maneuvers with precision and grace.
Off the mat, Mr. Bubba remains an enigma, rarely speaking or revealing his true identity. This air of mystery only adds to his allure, leaving fans and fellow wrestlers alike wondering what secrets lie beneath the bunny mask.
Whether he's captivating audiences with his mesmerizing performances or outwitting opponents with his unconventional style, Mr. Bubba is a force to be reckoned with in the world of professional wrestling. With his unique blend of creativity and athleticism, he is sure to leave a lasting impression on the industry for years to come."

# --- Skills: skill ratings for Mr Bubba, row 5 ---
$wsSkills.Cells.Item(5, 1).Value = 2964
$wsSkills.Cells.Item(5, 2).Value = 53
$wsSkills.Cells.Item(5, 3).Value = 32
$wsSkills.Cells.Item(5, 4).Value = 22
$wsSkills.Cells.Item(5, 5).Value = 56
$wsSkills.Cells.Item(5, 6).Value = 34
$wsSkills.Cells.Item(5, 7).Value = 65
$wsSkills.Cells.Item(5, 8).Value = 34
$wsSkills.Cells.Item(5, 9).Value = 82
$wsSkills.Cells.Item(5, 10).Value = 37
$wsSkills.Cells.Item(5, 11).Value = 53
$wsSkills.Cells.Item(5, 12).Value = 89
$wsSkills.Cells.Item(5, 13).Value = 81
$wsSkills.Cells.Item(5, 14).Value = 33
$wsSkills.Cells.Item(5, 15).Value = 100
$wsSkills.Cells.Item(5, 16).Value = 100
$wsSkills.Cells.Item(5, 17).Value = 65
$wsSkills.Cells.Item(5, 18).Value = 78
$wsSkills.Cells.Item(5, 19).Value = 81
$wsSkills.Cells.Item(5, 20).Value = 33
$wsSkills.Cells.Item(5, 21).Value = 68
$wsSkills.Cells.Item(5, 22).Value = 51
$wsSkills.Cells.Item(5, 23).Value = 55
$wsSkills.Cells.Item(5, 24).Value = 81
$wsSkills.Cells.Item(5, 25).Value = 60
$wsSkills.Cells.Item(5, 26).Value = 0
$wsSkills.Cells.Item(5, 27).Value = 0
$wsSkills.Cells.Item(5, 28).Value = 0
$wsSkills.Cells.Item(5, 29).Value = 100
$wsSkills.Cells.Item(5, 30).Value = 0
$wsSkills.Cells.Item(5, 31).Value = 0
$wsSkills.Cells.Item(5, 32).Value = 0
$wsSkills.Cells.Item(5, 33).Value = 0
$wsSkills.Cells.Item(5, 34).Value = 0
$wsSkills.Cells.Item(5, 35).Value = 0
$wsSkills.Cells.Item(5, 36).Value = 0
$wsSkills.Cells.Item(5, 37).Value = 0
$wsSkills.Cells.Item(5, 38).Value = 6
$wsSkills.Cells.Item(5, 39).Value = 6
$wsSkills.Cells.Item(5, 40).Value = 6
$wsSkills.Cells.Item(5, 41).Value = 6
$wsSkills.Cells.Item(5, 42).Value = 6

# --- Contracts: fix existing row 4 (Tiger Mask) ContractDebutDate cell format ---
$wsContracts.Cells.Item(4, 31).Value = "00:00:00"
$wsContracts.Cells.Item(4, 31).ClearFormats()

# --- Contracts: new contract for Mr Bubba, row 5 ---
$wsContracts.Cells.Item(5, 1).Value = 2653
$wsContracts.Cells.Item(5, 2).Value = 119
$wsContracts.Cells.Item(5, 3).Value = 2964
$wsContracts.Cells.Item(5, 4).Value = "Mr Bubba"
$wsContracts.Cells.Item(5, 5).Value = "Mr"
$wsContracts.Cells.Item(5, 6).Value = "mrbubba.jpg"
$wsContracts.Cells.Item(5, 7).Value = 1
$wsContracts.Cells.Item(5, 8).Value = $false
$wsContracts.Cells.Item(5, 9).Value = 0
$wsContracts.Cells.Item(5, 10).Value = 0
$wsContracts.Cells.Item(5, 11).Value = 0
$wsContracts.Cells.Item(5, 12).Value = $true
$wsContracts.Cells.Item(5, 13).Value = $false
$wsContracts.Cells.Item(5, 14).Value = $false
$wsContracts.Cells.Item(5, 15).Value = $true
$wsContracts.Cells.Item(5, 16).Value = $false
$wsContracts.Cells.Item(5, 17).Value = $false
$wsContracts.Cells.Item(5, 18).Value = 0
$wsContracts.Cells.Item(5, 19).Value = 0
$wsContracts.Cells.Item(5, 20).Value = 0
$wsContracts.Cells.Item(5, 21).Value = 0
$wsContracts.Cells.Item(5, 22).Value = 0
$wsContracts.Cells.Item(5, 23).Value = $false
$wsContracts.Cells.Item(5, 24).Value = $false
$wsContracts.Cells.Item(5, 25).Value = $false
$wsContracts.Cells.Item(5, 26).Value = $false
$wsContracts.Cells.Item(5, 27).Value = 43446
$wsContracts.Cells.Item(5, 27).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsContracts.Cells.Item(5, 28).Value = 253
$wsContracts.Cells.Item(5, 29).Value = 0
$wsContracts.Cells.Item(5, 30).Value = 0
$wsContracts.Cells.Item(5, 31).Value = 0
$wsContracts.Cells.Item(5, 31).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsContracts.Cells.Item(5, 32).Value = -1
$wsContracts.Cells.Item(5, 33).Value = -1
$wsContracts.Cells.Item(5, 34).Value = 0
$wsContracts.Cells.Item(5, 35).Value = 0
$wsContracts.Cells.Item(5, 36).Value = 5
$wsContracts.Cells.Item(5, 37).Value = 0
$wsContracts.Cells.Item(5, 38).Value = $true
$wsContracts.Cells.Item(5, 39).Value = $false
$wsContracts.Cells.Item(5, 40).Value = $false
$wsContracts.Cells.Item(5, 41).Value = $false
$wsContracts.Cells.Item(5, 42).Value = $false
$wsContracts.Cells.Item(5, 43).Value = $false
$wsContracts.Cells.Item(5, 44).Value = $false
$wsContracts.Cells.Item(5, 45).Value = $false
$wsContracts.Cells.Item(5, 46).Value = $false
$wsContracts.Cells.Item(5, 47).Value = 200
$wsContracts.Cells.Item(5, 49).Value = 4
$wsContracts.Cells.Item(5, 50).Value = 2
$wsContracts.Cells.Item(5, 51).Value = 0
$wsContracts.Cells.Item(5, 52).Value = 0
$wsContracts.Cells.Item(5, 53).Value = 0
$wsContracts.Cells.Item(5, 54).Value = 0
$wsContracts.Cells.Item(5, 55).Value = 0
$wsContracts.Cells.Item(5, 56).Value = 0
$wsContracts.Cells.Item(5, 57).Value = $false
$wsContracts.Cells.Item(5, 58).Value = $false
$wsContracts.Cells.Item(5, 59).Value = $false
$wsContracts.Cells.Item(5, 60).Value = $false
$wsContracts.Cells.Item(5, 61).Value = $false
$wsContracts.Cells.Item(5, 62).Value = $false
$wsContracts.Cells.Item(5, 63).Value = $false
$wsContracts.Cells.Item(5, 64).Value = $false
$wsContracts.Cells.Item(5, 65).Value = $false
$wsContracts.Cells.Item(5, 66).Value = $false
$wsContracts.Cells.Item(5, 67).Value = $false
$wsContracts.Cells.Item(5, 68).Value = $false
$wsContracts.Cells.Item(5, 69).Value = $false
$wsContracts.Cells.Item(5, 70).Value = $false
$wsContracts.Cells.Item(5, 71).Value = $false
$wsContracts.Cells.Item(5, 72).Value = $false
$wsContracts.Cells.Item(5, 73).Value = $false
$wsContracts.Cells.Item(5, 74).Value = $false
$wsContracts.Cells.Item(5, 75).Value = $false
$wsContracts.Cells.Item(5, 76).Value = $false
$wsContracts.Cells.Item(5, 77).Value = $false
$wsContracts.Cells.Item(5, 78).Value = $false
$wsContracts.Cells.Item(5, 79).Value = $false
$wsContracts.Cells.Item(5, 80).Value = $false
$wsContracts.Cells.Item(5, 81).Value = $false

# --- Notes: preset/gimmick info for Mr Bubba, row 5 ---
$wsNotes.Cells.Item(5, 1).Value = "Mr Bubba"
$wsNotes.Cells.Item(5, 2).Value = "A new kind of wrestler that wears a bunny mask"
$wsNotes.Cells.Item(5, 3).Value = "Male"
$wsNotes.Cells.Item(5, 4).Value = "Bubba"
$wsNotes.Cells.Item(5, 5).Value = "Random"
$wsNotes.Cells.Item(5, 6).Value = "Interpret"
$wsNotes.Cells.Item(5, 7).Value = "mrbubba.jpg"
$wsNotes.Cells.Item(5, 8).Value = "Mr Bubba stands at 6 feet tall, with a muscular build and a mysterious aura emanating from his intimidating stature and unique bunny mask."
$wsNotes.Cells.Item(5, 9).Value = $true
$wsNotes.Cells.Item(5, 10).Value = 9


